$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set updated / new cell values (diagonal shift introducing ifoCAST full series column)
$ws.Range("K2").Value = 0.1277855347458341
$ws.Range("B3").Value = 9.643547872324863
$ws.Range("C3").Value = -8.527713777675137
$ws.Range("D3").Value = -0.3880427776751375
$ws.Range("E3").Value = 1.759528090324862
$ws.Range("F3").Value = -1.680501777675137
$ws.Range("G3").Value = -1.741675777675137
$ws.Range("H3").Value = 0.2997802223248625
$ws.Range("I3").Value = -0.2828147776751375
$ws.Range("J3").Value = -0.1871396336751375
$ws.Range("K3").Value = -0.4475587776751375
$ws.Range("B4").Value = -9.584088737857668
$ws.Range("C4").Value = -1.444417737857667
$ws.Range("D4").Value = 0.7031531301423328
$ws.Range("E4").Value = -2.736876737857667
$ws.Range("F4").Value = -2.798050737857667
$ws.Range("G4").Value = -0.7565947378576671
$ws.Range("H4").Value = -1.339189737857667
$ws.Range("I4").Value = -1.243514593857667
$ws.Range("J4").Value = -1.503933737857667
$ws.Range("K4").Value = -0.6684787378576671
$ws.Range("B5").Value = -0.5737713618867365
$ws.Range("C5").Value = 1.573799506113263
$ws.Range("D5").Value = -1.866230361886736
$ws.Range("E5").Value = -1.927404361886736
$ws.Range("F5").Value = 0.1140516381132635
$ws.Range("G5").Value = -0.4685433618867365
$ws.Range("H5").Value = -0.3728682178867365
$ws.Range("I5").Value = -0.6332873618867365
$ws.Range("J5").Value = 0.2021676381132635
$ws.Range("K5").Value = 0.1020946851132635
$ws.Range("B6").Value = 1.56314828978338
$ws.Range("C6").Value = -1.87688157821662
$ws.Range("D6").Value = -1.93805557821662
$ws.Range("E6").Value = 0.10340042178338
$ws.Range("F6").Value = -0.47919457821662
$ws.Range("G6").Value = -0.38351943421662
$ws.Range("H6").Value = -0.64393857821662
$ws.Range("I6").Value = 0.19151642178338
$ws.Range("J6").Value = 0.09144346878338
$ws.Range("K6").Value = -0.26234957821662
$ws.Range("B7").Value = -1.985496110432761
$ws.Range("C7").Value = -2.046670110432761
$ws.Range("D7").Value = -0.005214110432760821
$ws.Range("E7").Value = -0.5878091104327609
$ws.Range("F7").Value = -0.4921339664327609
$ws.Range("G7").Value = -0.7525531104327607
$ws.Range("H7").Value = 0.08290188956723921
$ws.Range("I7").Value = -0.01717106343276081
$ws.Range("J7").Value = -0.3709641104327608
$ws.Range("K7").Value = -0.2227861104327608
$ws.Range("B8").Value = -1.838568389465022
$ws.Range("C8").Value = 0.2028876105349782
$ws.Range("D8").Value = -0.3797073894650218
$ws.Range("E8").Value = -0.2840322454650218
$ws.Range("F8").Value = -0.5444513894650218
$ws.Range("G8").Value = 0.2910036105349783
$ws.Range("H8").Value = 0.1909306575349782
$ws.Range("I8").Value = -0.1628623894650218
$ws.Range("J8").Value = -0.01468438946502179
$ws.Range("K8").Value = 0.1427616105349782
$ws.Range("B9").Value = 0.3193857809087786
$ws.Range("C9").Value = -0.2632092190912214
$ws.Range("D9").Value = -0.1675340750912214
$ws.Range("E9").Value = -0.4279532190912214
$ws.Range("F9").Value = 0.4075017809087786
$ws.Range("G9").Value = 0.3074288279087786
$ws.Range("H9").Value = -0.04636421909122137
$ws.Range("I9").Value = 0.1018137809087786
$ws.Range("J9").Value = 0.2592597809087786
$ws.Range("K9").Value = -0.2414032190912214
$ws.Range("B10").Value = -0.2177158050940863
$ws.Range("C10").Value = -0.1220406610940863
$ws.Range("D10").Value = -0.3824598050940862
$ws.Range("E10").Value = 0.4529951949059138
$ws.Range("F10").Value = 0.3529222419059138
$ws.Range("G10").Value = -0.0008708050940862688
$ws.Range("H10").Value = 0.1473071949059137
$ws.Range("I10").Value = 0.3047531949059137
$ws.Range("J10").Value = -0.1959098050940863
$ws.Range("K10").Value = 0.08445119490591374
$ws.Range("B11").Value = -0.1395947824540507
$ws.Range("C11").Value = -0.4000139264540507
$ws.Range("D11").Value = 0.4354410735459493
$ws.Range("E11").Value = 0.3353681205459493
$ws.Range("F11").Value = -0.01842492645405072
$ws.Range("G11").Value = 0.1297530735459493
$ws.Range("H11").Value = 0.2871990735459493
$ws.Range("I11").Value = -0.2134639264540507
$ws.Range("J11").Value = 0.06689707354594927
$ws.Range("K11").Value = -0.1050779264540507
$ws.Range("B12").Value = -0.3119066075030444
$ws.Range("C12").Value = 0.5235483924969555
$ws.Range("D12").Value = 0.4234754394969555
$ws.Range("E12").Value = 0.06968239249695551
$ws.Range("F12").Value = 0.2178603924969555
$ws.Range("G12").Value = 0.3753063924969555
$ws.Range("H12").Value = -0.1253566075030445
$ws.Range("I12").Value = 0.1550043924969555
$ws.Range("J12").Value = -0.01697060750304449
$ws.Range("K12").Value = 0.2888923924969555
$ws.Range("B13").Value = 0.7021230259847391
$ws.Range("C13").Value = 0.6020500729847391
$ws.Range("D13").Value = 0.2482570259847391
$ws.Range("E13").Value = 0.3964350259847391
$ws.Range("F13").Value = 0.5538810259847391
$ws.Range("G13").Value = 0.0532180259847391
$ws.Range("H13").Value = 0.3335790259847391
$ws.Range("I13").Value = 0.1616040259847391
$ws.Range("J13").Value = 0.4674670259847391
$ws.Range("K13").Value = -0.1448629740152609
$ws.Range("B14").Value = 1.514077300737389
$ws.Range("C14").Value = 1.160284253737389
$ws.Range("D14").Value = 1.308462253737389
$ws.Range("E14").Value = 1.465908253737389
$ws.Range("F14").Value = 0.9652452537373891
$ws.Range("G14").Value = 1.245606253737389
$ws.Range("H14").Value = 1.073631253737389
$ws.Range("I14").Value = 1.379494253737389
$ws.Range("J14").Value = 0.7671642537373891
$ws.Range("K14").Value = 1.455535253737389
$ws.Range("B15").Value = 0.2163100177716323
$ws.Range("C15").Value = 0.3644880177716323
$ws.Range("D15").Value = 0.5219340177716323
$ws.Range("E15").Value = 0.02127101777163229
$ws.Range("F15").Value = 0.3016320177716323
$ws.Range("G15").Value = 0.1296570177716323
$ws.Range("H15").Value = 0.4355200177716323
$ws.Range("I15").Value = -0.1768099822283677
$ws.Range("J15").Value = 0.5115610177716323
$ws.Range("B16").Value = 0.3684559152847414
$ws.Range("C16").Value = 0.5259019152847414
$ws.Range("D16").Value = 0.0252389152847414
$ws.Range("E16").Value = 0.3055999152847414
$ws.Range("F16").Value = 0.1336249152847414
$ws.Range("G16").Value = 0.4394879152847414
$ws.Range("H16").Value = -0.1728420847152586
$ws.Range("I16").Value = 0.5155289152847414
$ws.Range("B17").Value = 0.6615420054549828
$ws.Range("C17").Value = 0.1608790054549828
$ws.Range("D17").Value = 0.4412400054549828
$ws.Range("E17").Value = 0.2692650054549828
$ws.Range("F17").Value = 0.5751280054549828
$ws.Range("G17").Value = -0.03720199454501721
$ws.Range("H17").Value = 0.6511690054549828
$ws.Range("B18").Value = -0.07992717640068192
$ws.Range("C18").Value = 0.2004338235993181
$ws.Range("D18").Value = 0.02845882359931809
$ws.Range("E18").Value = 0.3343218235993181
$ws.Range("F18").Value = -0.2780081764006819
$ws.Range("G18").Value = 0.4103628235993181
$ws.Range("B19").Value = 0.1551026083472217
$ws.Range("C19").Value = -0.01687239165277828
$ws.Range("D19").Value = 0.2889906083472217
$ws.Range("E19").Value = -0.3233393916527783
$ws.Range("F19").Value = 0.3650316083472218
$ws.Range("B20").Value = -0.08373381358883215
$ws.Range("C20").Value = 0.2221291864111679
$ws.Range("D20").Value = -0.3902008135888321
$ws.Range("E20").Value = 0.2981701864111679
$ws.Range("B21").Value = 0.1925429840544512
$ws.Range("C21").Value = -0.4197870159455487
$ws.Range("D21").Value = 0.2685839840544513
$ws.Range("B22").Value = -0.4379376726751896
$ws.Range("C22").Value = 0.2504333273248104
$ws.Range("B23").Value = 0.2324015030761111

# Clear cells that are no longer present (row shrinks by one column at the tail)
$ws.Range("K15").ClearContents()
$ws.Range("J16").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("H18").ClearContents()
$ws.Range("G19").ClearContents()
$ws.Range("F20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("B24").ClearContents()

Write-Host "Applied ifoCAST full series evaluation update"
